$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 449; existing rows 449-501 shift down to 450-502,
# preserving all of their data and formatting.
$ws.Rows.Item(449).Insert()

# Populate the newly inserted row 449 with the new record.
$ws.Range("A449").Value = 10
$ws.Range("B449").Value = "Vega Modelo de Temuco"
$ws.Range("C449").Value = "La Araucanía"
$ws.Range("D449").Value = 45124
$ws.Range("E449").Value = 9
$ws.Range("F449").Value = 100114013
$ws.Range("G449").Value = "Zanahoria"
$ws.Range("H449").Value = "Sin especificar"
$ws.Range("I449").Value = "Primera"
$ws.Range("J449").Value = 150
$ws.Range("K449").Value = 5000
$ws.Range("L449").Value = 5000
$ws.Range("M449").Value = 5000
$ws.Range("N449").Value = "$/saco 25 kilos"
$ws.Range("O449").Value = "Región de La Araucanía"
$ws.Range("P449").Value = 200
$ws.Range("Q449").Value = 25
$ws.Range("R449").Value = "Hortaliza"
